$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 480, shifting rows 480:492 down to 481:493
$ws.Rows.Item(480).Insert()

# Populate the newly inserted row 480 with the new weekly record
$ws.Cells.Item(480, 1).Value = 4
$ws.Cells.Item(480, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(480, 3).Value = "Los Lagos"
$ws.Cells.Item(480, 4).Value = 45239
$ws.Cells.Item(480, 4).Style = $ws.Cells.Item(481, 4).Style
$ws.Cells.Item(480, 5).Value = 10
$ws.Cells.Item(480, 6).Value = "Fruta"
$ws.Cells.Item(480, 7).Value = 100108
$ws.Cells.Item(480, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(480, 9).Value = 100108005
$ws.Cells.Item(480, 10).Value = "Piña"
$ws.Cells.Item(480, 11).Value = "Caramelo"
$ws.Cells.Item(480, 12).Value = "Primera"
$ws.Cells.Item(480, 13).Value = 100
$ws.Cells.Item(480, 14).Value = 25000
$ws.Cells.Item(480, 15).Value = 25000
$ws.Cells.Item(480, 16).Value = 25000
$ws.Cells.Item(480, 17).Value = "$/caja 12 unidades"
$ws.Cells.Item(480, 18).Value = "Ecuador"
$ws.Cells.Item(480, 19).Value = 2083
$ws.Cells.Item(480, 20).Value = 12
